# Applies the scheduled update of the date and the practice problems
# for the "three-digit number multiplied by one-digit number" worksheet.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-01-30 Thursday"; New = "2025-01-31 Friday" },
    @{ Old = "892×7="; New = "456×7=" },
    @{ Old = "319×2="; New = "381×9=" },
    @{ Old = "184×7="; New = "635×8=" },
    @{ Old = "189×2="; New = "247×9=" },
    @{ Old = "321×6="; New = "291×9=" },
    @{ Old = "978×6="; New = "265×4=" },
    @{ Old = "107×6="; New = "644×7=" },
    @{ Old = "866×7="; New = "841×4=" },
    @{ Old = "686×6="; New = "196×3=" },
    @{ Old = "264×9="; New = "314×3=" },
    @{ Old = "162×3="; New = "248×3=" },
    @{ Old = "472×5="; New = "126×4=" },
    @{ Old = "340×5="; New = "207×5=" },
    @{ Old = "101×9="; New = "225×4=" },
    @{ Old = "114×8="; New = "178×5=" },
    @{ Old = "793×7="; New = "429×6=" },
    @{ Old = "161×3="; New = "869×2=" },
    @{ Old = "520×2="; New = "629×4=" },
    @{ Old = "684×2="; New = "654×8=" },
    @{ Old = "593×9="; New = "124×7=" },
    @{ Old = "296×5="; New = "620×2=" },
    @{ Old = "859×2="; New = "744×9=" },
    @{ Old = "172×3="; New = "886×3=" },
    @{ Old = "566×3="; New = "813×2=" },
    @{ Old = "243×4="; New = "108×2=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
